$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive leading text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Working on waypoints*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$r = $target.Range

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">I changed the </w:t></w:r><w:r><w:t>read</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>from</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve">file code to read </w:t></w:r><w:r><w:t xml:space="preserve">the patrol </w:t></w:r><w:r><w:t xml:space="preserve">waypoint box </w:t></w:r><w:r><w:t>numbers</w:t></w:r><w:r><w:t xml:space="preserve"> from </w:t></w:r><w:r><w:t xml:space="preserve">the map </w:t></w:r><w:r><w:t xml:space="preserve">file and store </w:t></w:r><w:r><w:t>the corresponding boxes in a list of waypoints in the box world class</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">While doing so, I </w:t></w:r><w:r><w:t xml:space="preserve">got sick of the low framerate, </w:t></w:r><w:r><w:t>I changed box.draw()</w:t></w:r><w:r><w:t xml:space="preserve"> to not fill the circle</w:t></w:r><w:r><w:t xml:space="preserve"> rendered for walls</w:t></w:r><w:r><w:t>, but just do a thick outline, saving iterations of circle outline drawing</w:t></w:r><w:r><w:t xml:space="preserve"> and increasing the framerate</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>I altered the agent setup code so that the soldier leader would have its target set at the start and plan a path to the first waypoint, prompting the rest of the soldiers to follow it. Then I reorganised the various conditions for planning a new path into one method that could manage target selection and path planning for all of its preceding code, and replaced them with a call to that method.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>This way, the target selection and call of self.plan_path() is all in one location and can be more easily modified as needed.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$r.InsertXML($xml)
